# Regenerate the "K" column (column G) values in the save_data sheet.
# These are the true strikeout ("K") counts that replace the previous
# "Strike#" placeholder figures used in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value for column G (rows 2-32)
$kValues = [ordered]@{
    2  = 5
    3  = 3
    4  = 3
    5  = 6
    6  = 3
    7  = 2
    8  = 3
    9  = 1
    10 = 7
    11 = 1
    12 = 8
    13 = 3
    14 = 8
    15 = 7
    16 = 9
    17 = 0
    18 = 1
    19 = 3
    20 = 4
    21 = 4
    22 = 5
    23 = 1
    24 = 4
    25 = 0
    26 = 3
    27 = 3
    28 = 5
    29 = 1
    30 = 4
    31 = 5
    32 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
